$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This edit updates the weekly price records for rows 45-64 on the
# "Arveja Verde" sheet. A new record was inserted at the top of this
# date-ordered block (row 45), shifting the previously-recorded values
# for rows 45-63 down into rows 46-64 (row 64's original record is
# superseded / dropped, matching the published diff).

# Row 45
$ws.Range("D45").Value = 44784
$ws.Range("J45").Value = 360
$ws.Range("K45").Value = 27000
$ws.Range("L45").Value = 29000
$ws.Range("M45").Value = 28000
$ws.Range("P45").Value = 1120

# Row 46
$ws.Range("D46").Value = 44685
$ws.Range("J46").Value = 160
$ws.Range("K46").Value = 25000
$ws.Range("L46").Value = 27000
$ws.Range("M46").Value = 26000
$ws.Range("P46").Value = 1040

# Row 47
$ws.Range("D47").Value = 44392
$ws.Range("J47").Value = 100
$ws.Range("K47").Value = 26000
$ws.Range("L47").Value = 28000
$ws.Range("M47").Value = 27000
$ws.Range("P47").Value = 1080

# Row 48
$ws.Range("D48").Value = 44370
$ws.Range("J48").Value = 400
$ws.Range("K48").Value = 27000
$ws.Range("L48").Value = 28000
$ws.Range("M48").Value = 27500
$ws.Range("P48").Value = 1100

# Row 49
$ws.Range("D49").Value = 44384
$ws.Range("J49").Value = 400
$ws.Range("K49").Value = 26000
$ws.Range("L49").Value = 28000
$ws.Range("M49").Value = 27000
$ws.Range("P49").Value = 1080

# Row 50
$ws.Range("D50").Value = 44748
$ws.Range("J50").Value = 700
$ws.Range("K50").Value = 28000
$ws.Range("L50").Value = 30000
$ws.Range("M50").Value = 29000
$ws.Range("P50").Value = 1160

# Row 51
$ws.Range("D51").Value = 44371
$ws.Range("J51").Value = 500
$ws.Range("K51").Value = 28000
$ws.Range("L51").Value = 30000
$ws.Range("M51").Value = 29000
$ws.Range("P51").Value = 1160

# Row 52
$ws.Range("D52").Value = 44427
$ws.Range("J52").Value = 300
$ws.Range("K52").Value = 28000
$ws.Range("L52").Value = 30000
$ws.Range("M52").Value = 29000
$ws.Range("P52").Value = 1160

# Row 53
$ws.Range("D53").Value = 44441
$ws.Range("J53").Value = 700
$ws.Range("K53").Value = 28000
$ws.Range("L53").Value = 30000
$ws.Range("M53").Value = 29000
$ws.Range("P53").Value = 1160

# Row 54
$ws.Range("D54").Value = 44419
$ws.Range("J54").Value = 600
$ws.Range("K54").Value = 27000
$ws.Range("L54").Value = 29000
$ws.Range("M54").Value = 28000
$ws.Range("P54").Value = 1120

# Row 55
$ws.Range("D55").Value = 44412
$ws.Range("J55").Value = 600
$ws.Range("K55").Value = 25000
$ws.Range("L55").Value = 27000
$ws.Range("M55").Value = 26000
$ws.Range("P55").Value = 1040

# Row 56
$ws.Range("D56").Value = 44483
$ws.Range("J56").Value = 300
$ws.Range("K56").Value = 18000
$ws.Range("L56").Value = 20000
$ws.Range("M56").Value = 19000
$ws.Range("P56").Value = 760

# Row 57
$ws.Range("D57").Value = 44469
$ws.Range("J57").Value = 600
$ws.Range("K57").Value = 22000
$ws.Range("L57").Value = 24000
$ws.Range("M57").Value = 23000
$ws.Range("P57").Value = 920

# Row 58
$ws.Range("D58").Value = 44434
$ws.Range("J58").Value = 500
$ws.Range("K58").Value = 28000
$ws.Range("L58").Value = 30000
$ws.Range("M58").Value = 29000
$ws.Range("P58").Value = 1160

# Row 59
$ws.Range("D59").Value = 44776
$ws.Range("J59").Value = 400
$ws.Range("K59").Value = 28000
$ws.Range("L59").Value = 30000
$ws.Range("M59").Value = 29000
$ws.Range("P59").Value = 1160

# Row 60
$ws.Range("D60").Value = 44356
$ws.Range("J60").Value = 300
$ws.Range("K60").Value = 26000
$ws.Range("L60").Value = 28000
$ws.Range("M60").Value = 27000
$ws.Range("P60").Value = 1080

# Row 61
$ws.Range("D61").Value = 44769
$ws.Range("J61").Value = 500
$ws.Range("K61").Value = 30000
$ws.Range("L61").Value = 32000
$ws.Range("M61").Value = 31000
$ws.Range("P61").Value = 1240

# Row 62
$ws.Range("D62").Value = 44399
$ws.Range("J62").Value = 400
$ws.Range("K62").Value = 26000
$ws.Range("L62").Value = 28000
$ws.Range("M62").Value = 27000
$ws.Range("P62").Value = 1080

# Row 63
$ws.Range("D63").Value = 44783
$ws.Range("J63").Value = 400
$ws.Range("K63").Value = 27000
$ws.Range("L63").Value = 29000
$ws.Range("M63").Value = 28000
$ws.Range("P63").Value = 1120

# Row 64
$ws.Range("D64").Value = 44377
$ws.Range("J64").Value = 500
$ws.Range("K64").Value = 26000
$ws.Range("L64").Value = 28000
$ws.Range("M64").Value = 27000
$ws.Range("P64").Value = 1080
